$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Append a trailing space run to the "For webapp team..." paragraph
# ------------------------------------------------------------------
$pWebapp = $d.Paragraphs.Item(27)
$rWebapp = $pWebapp.Range
$rWebapp.InsertAfter(" ")

# ------------------------------------------------------------------
# 2. Insert four new sub-bullets (ilvl=2) right after that paragraph,
#    before "For simulation team..."
# ------------------------------------------------------------------
$ip = $rWebapp.InsertParagraphAfter()
$pA = $d.Paragraphs.Item(28)
$pA.Range.ListFormat.ListIndent()
$pA.Range.Text = "Designs may have more time than the simulation, entirely depending on Laskey’s (and possibly Garfield’s) availability."

$ip = $pA.Range.InsertParagraphAfter()
$pB = $d.Paragraphs.Item(29)
$pB.Range.Text = "These will be presented to Laskey as options for her. At least 3 designs, ideally 5, but I understand that it is difficult to come up with ideas"

$ip = $pB.Range.InsertParagraphAfter()
$pC = $d.Paragraphs.Item(30)
$pC.Range.Text = "Understand that these designs will not be concrete, they will be options on what the webapp will basically look like."
$rItalic = $pC.Range.Duplicate
$rItalic.Find.Execute("basically", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rItalic.Font.Italic = 1

$ip = $pC.Range.InsertParagraphAfter()
$pD = $d.Paragraphs.Item(31)
$pD.Range.Text = "Maybe add color schemes, but this isn’t as important as what a design will look like "
# Word recorded a page break inside this run while paginating the
# now-longer document; reproduce that bookkeeping element.
$pDFull = $pD.Range
$pDXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Maybe add color schemes, but this isn' + [char]0x2019 + 't as important as what a design will look like </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pDFull.InsertXML($pDXml)

# ------------------------------------------------------------------
# 3. Insert seven new bullets after "For simulation team..."
# ------------------------------------------------------------------
$pSim = $d.Paragraphs.Item(32)

$ip = $pSim.Range.InsertParagraphAfter()
$pE = $d.Paragraphs.Item(33)
$pE.Range.ListFormat.ListIndent()
$pE.Range.Text = "Only looking for a high level overview by next Tuesday (10/11/2022)"

$ip = $pE.Range.InsertParagraphAfter()
$pF = $d.Paragraphs.Item(34)
$pF.Range.ListFormat.ListIndent()
$pF.Range.Text = "What systems talk to which elements."

$ip = $pF.Range.InsertParagraphAfter()
$pG = $d.Paragraphs.Item(35)
$pG.Range.Text = "How users interact with the system."

$ip = $pG.Range.InsertParagraphAfter()
$pH = $d.Paragraphs.Item(36)
$pH.Range.Text = "How elements interact with each other."

$ip = $pH.Range.InsertParagraphAfter()
$pI = $d.Paragraphs.Item(37)
$pI.Range.ListFormat.ListOutdent()
$pI.Range.ListFormat.ListOutdent()
$pI.Range.ListFormat.ListOutdent()
$pI.Range.Text = "Goal is to be programming the simulation by next week (week of 10/10/2022), and hopefully start the website, but that is more on the availability of Dr Laskey."

$ip = $pI.Range.InsertParagraphAfter()
$pJ = $d.Paragraphs.Item(38)
$pJ.Range.Text = "IF WE DO NOT HEAR FROM LASKEY OR GARFIELD REGARDING THE DESIGN, AS A TEAM WE WILL VOTE ON THE DESIGN WE LIKE THE MOST AND THAT WILL BE THE DESIGN WE WILL GO WITH"

$ip = $pJ.Range.InsertParagraphAfter()
$pK = $d.Paragraphs.Item(39)
$pK.Range.ListFormat.ListIndent()
$pK.Range.Text = "We are not starting the website any later than 10/18/2022."
